# Planejamento.xlsx update
# "refatoracao do codigo, complemento ao movimento do sprites, pontuacao, vidas"
#
# Updates the status-tracking grid (columns C/D/E = Anne/Matheus/Juliane) for
# rows 11-14 so completed items are marked with "x" (rendered in the default
# black font) instead of "A FAZER" (rendered in the red/theme font), and
# reassigns a couple of the remaining "A FAZER" / "x" markers between
# columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - "Posicionar aleatoriamente x fantasminhas e y pac's ao iniciar o jogo"
$ws.Range("C11").Value = "x"
$ws.Range("C11").Font.Color = 0
$ws.Range("D11").Font.Color = 0
$ws.Range("E11").Value = "x"
$ws.Range("E11").Font.Color = 0

# Row 12 - "Fazer os movimentos dos fantasminhas e dos pac's: horizontal, vertical, diagonal"
$ws.Range("C12").Value = "x"
$ws.Range("C12").Font.Color = 0
$ws.Range("D12").Font.Color = 0
$ws.Range("E12").Value = "x"
$ws.Range("E12").Font.Color = 0

# Row 13 - "Dar sentido ao tap: fantasminha = explosao; pac = brilho"
$ws.Range("C13").Font.Color = 0
$ws.Range("D13").Value = "x"
$ws.Range("D13").Font.Color = 0
$ws.Range("E13").Font.Color = 0

# Row 14 - "Trabalhar nos extras da tela de jogo: pontuacao, vida, cronometro"
$ws.Range("D14").ClearContents() | Out-Null
$ws.Range("E14").Value = "x"
$ws.Range("E14").Font.Color = 0
